$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: dependent coordinate system template (cs fields blank, dependence = dependent, co_x/co_y = false)
$ws.Range("A2").Value = "'"
$ws.Range("B2").Value = "'"
$ws.Range("C2").Value = "dependent"
$ws.Range("D2").Value = "'"
$ws.Range("E2").Value = "'"
$ws.Range("F2").Value = "'"
$ws.Range("G2").Value = "'false"
$ws.Range("H2").Value = "'false"
$ws.Range("A2:H2").Style = "Normal"

# Row 3: independent coordinate system template (cs fields blank, dependence = independent, co_x/co_y = true)
$ws.Range("A3").Value = "'"
$ws.Range("B3").Value = "'"
$ws.Range("C3").Value = "independent"
$ws.Range("D3").Value = "'"
$ws.Range("E3").Value = "'"
$ws.Range("F3").Value = "'"
$ws.Range("G3").Value = "'true"
$ws.Range("H3").Value = "'true"
$ws.Range("A3:H3").Style = "Normal"
